# Updated cryptos list values (price + volume%) and a few coin row swaps/renames,
# reproducing the scraped-data refresh described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.812.95'
$ws.Range('E2').Value = '  -2.55%  '

$ws.Range('D3').Value = '2.748.74'
$ws.Range('E3').Value = '  -1.69%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').Value = "'" + '349.58'
$ws.Range('E5').Value = '  -3.34%  '

$ws.Range('D6').Value = "'" + '106.66'
$ws.Range('E6').Value = '  -3.08%  '

$ws.Range('D7').Value = "'" + '0.545'
$ws.Range('E7').Value = '  -3.12%  '

$ws.Range('D8').Value = "'" + '1.00'
$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').Value = "'" + '0.578'
$ws.Range('E9').Value = '  -3.17%  '

$ws.Range('D10').Value = "'" + '39.08'
$ws.Range('E10').Value = '  -3.12%  '

$ws.Range('E11').Value = '  +3.30%  '

$ws.Range('D12').Value = "'" + '0.0829'
$ws.Range('E12').Value = '  -3.09%  '

$ws.Range('D13').Value = "'" + '19.72'
$ws.Range('E13').Value = '  +0.66%  '

$ws.Range('D14').Value = "'" + '7.45'
$ws.Range('E14').Value = '  -3.10%  '

$ws.Range('D15').Value = '3.173.86'
$ws.Range('E15').Value = '  -2.04%  '

$ws.Range('D16').Value = '2.737.69'
$ws.Range('E16').Value = '  -2.71%  '

$ws.Range('D17').Value = "'" + '0.921'
$ws.Range('E17').Value = '  -1.79%  '

$ws.Range('D18').Value = '50.739.11'
$ws.Range('E18').Value = '  -2.57%  '

$ws.Range('D19').Value = "'" + '7.54'
$ws.Range('E19').Value = '  +1.83%  '

$ws.Range('E20').Value = '  -3.45%  '

$ws.Range('D21').Value = "'" + '12.93'
$ws.Range('E21').Value = '  -1.97%  '

$ws.Range('D22').Value = '0.0₃0953'
$ws.Range('E22').Value = '  -3.22%  '

$ws.Range('D23').Value = "'" + '69.26'
$ws.Range('E23').Value = '  -0.88%  '

$ws.Range('D24').Value = "'" + '263.24'
$ws.Range('E24').Value = '  -3.91%  '

$ws.Range('D25').Value = "'" + '2.69'
$ws.Range('E25').Value = '  -3.23%  '

$ws.Range('D26').Value = "'" + '1.00'
$ws.Range('E26').Value = '  +0.00%  '

$ws.Range('D27').Value = "'" + '25.78'
$ws.Range('E27').Value = '  -3.67%  '

$ws.Range('D28').Value = "'" + '0.161'
$ws.Range('E28').Value = '  +12.43%  '

$ws.Range('D29').Value = "'" + '10.05'
$ws.Range('E29').Value = '  -1.65%  '

$ws.Range('E30').Value = '  +0.24%  '

$ws.Range('D31').Value = "'" + '51.51'
$ws.Range('E31').Value = '  +0.00%  '

$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').Value = "'" + '34.16'
$ws.Range('E32').Value = '  -1.31%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'" + '5.97'
$ws.Range('E33').Value = '  +3.29%  '

$ws.Range('D34').Value = "'" + '0.0443'
$ws.Range('E34').Value = '  -6.45%  '

$ws.Range('E35').Value = '  -2.51%  '

$ws.Range('D36').Value = "'" + '0.0825'
$ws.Range('E36').Value = '  -2.27%  '

$ws.Range('E37').Value = '  -0.06%  '

$ws.Range('D38').Value = "'" + '18.21'
$ws.Range('E38').Value = '  -0.96%  '

$ws.Range('E39').Value = '  -3.45%  '

$ws.Range('E40').Value = '  -3.79%  '

$ws.Range('E41').Value = '  -1.59%  '

$ws.Range('D42').Value = "'" + '2.46'
$ws.Range('E42').Value = '  -4.88%  '

$ws.Range('D43').Value = "'" + '120.29'
$ws.Range('E43').Value = '  -3.59%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = "'" + '22.00'
$ws.Range('E44').Value = '  -0.55%  '

$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = "'" + '2.19'
$ws.Range('E45').Value = '  -2.37%  '

$ws.Range('D46').Value = '2.076.66'
$ws.Range('E46').Value = '  +0.24%  '

$ws.Range('D47').Value = "'" + '3.21'
$ws.Range('E47').Value = '  -1.99%  '

$ws.Range('E48').Value = '  -1.63%  '

$ws.Range('D49').Value = "'" + '0.906'
$ws.Range('E49').Value = '  -4.27%  '

$ws.Range('D50').Value = "'" + '5.41'
$ws.Range('E50').Value = '  -6.09%  '

$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = "'" + '57.33'
$ws.Range('E51').Value = '  -3.55%  '
